$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.228.78'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '1.663.15'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("ZZ1").Formula = "=TEXT(218.25,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.5229,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E6").Value = '  -1.59%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.2671,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.06296,""0.00000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("ZZ1").Formula = "=TEXT(20.83,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E10").Value = '  -4.49%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.07726,""0.00000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E11").Value = '  -1.19%  '
$ws.Range("D12").Value = '1.661.48'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("ZZ1").Formula = "=TEXT(4.437,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("D14").Value = '1.889.05'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.5454,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("D16").Value = '0.0₅8192'
$ws.Range("E16").Value = '  -1.81%  '
$ws.Range("ZZ1").Formula = "=TEXT(64.59,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").Value = '26.253.58'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("ZZ1").Formula = "=TEXT(4.651,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E20").Value = '  -1.92%  '
$ws.Range("ZZ1").Formula = "=TEXT(193.87,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("ZZ1").Formula = "=TEXT(10.08,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("ZZ1").Formula = "=TEXT(6.044,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E23").Value = '  -4.93%  '
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("ZZ1").Formula = "=TEXT(139.95,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.1232,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("ZZ1").Formula = "=TEXT(7.164,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E27").Value = '  -3.12%  '
$ws.Range("ZZ1").Formula = "=TEXT(16.14,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.06130,""0.00000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E30").Value = '  -3.04%  '
$ws.Range("ZZ1").Formula = "=TEXT(1.279,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("ZZ1").Formula = "=TEXT(3.266,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E33").Value = '  -5.47%  '
$ws.Range("ZZ1").Formula = "=TEXT(1.623,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E34").Value = '  -3.40%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.9681,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E35").Value = '  -4.28%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("ZZ1").Formula = "=TEXT(2.784,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.5683,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E38").Value = '  -8.68%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.01594,""0.00000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("ZZ1").Formula = "=TEXT(5.986,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.8568,""0.0000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E41").Value = '  -0.94%  '
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("D43").Value = '1.014.90'
$ws.Range("E43").Value = '  -6.99%  '
$ws.Range("ZZ1").Formula = "=TEXT(100.34,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").Value = '1.804.83'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("D46").Value = '0.0₈108'
$ws.Range("E46").Value = '  +5.21%  '
$ws.Range("ZZ1").Formula = "=TEXT(57.10,""0.00"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("ZZ1").Formula = "=TEXT(1.009,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E48").Value = '  +0.87%  '
$ws.Range("ZZ1").Formula = "=TEXT(7.979,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("ZZ1").Formula = "=TEXT(1.483,""0.000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("ZZ1").Formula = "=TEXT(0.05186,""0.00000"")"
$ws.Range("ZZ1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("E51").Value = '  -0.46%  '
